$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$data = @(
    @("PEOPLES R CHINA", 5089),
    @("USA", 1760),
    @("GERMANY", 550),
    @("INDIA", 550),
    @("SOUTH KOREA", 526),
    @("JAPAN", 393),
    @("ENGLAND", 359),
    @("NETHERLANDS", 288),
    @("IRAN", 286),
    @("CANADA", 280),
    @("BELGIUM", 241),
    @("AUSTRALIA", 234),
    @("FRANCE", 220),
    @("ITALY", 217),
    @("SPAIN", 175),
    @("SINGAPORE", 171),
    @("TAIWAN", 123),
    @("RUSSIA", 117),
    @("SWITZERLAND", 115),
    @("SAUDI ARABIA", 101),
    @("MALAYSIA", 93),
    @("TURKEY", 91),
    @("BRAZIL", 80),
    @("POLAND", 78),
    @("PORTUGAL", 64),
    @("SWEDEN", 63),
    @("THAILAND", 54),
    @("PAKISTAN", 51),
    @("EGYPT", 48),
    @("GREECE", 48),
    @("FINLAND", 42),
    @("ISRAEL", 40),
    @("AUSTRIA", 39),
    @("CZECH REPUBLIC", 38),
    @("DENMARK", 38),
    @("QATAR", 36),
    @("ROMANIA", 36),
    @("SCOTLAND", 34),
    @("TURKIYE", 30),
    @("ARGENTINA", 29),
    @("VIETNAM", 28),
    @("WALES", 25),
    @("U ARAB EMIRATES", 24),
    @("IRELAND", 22),
    @("NEW ZEALAND", 21),
    @("CHILE", 20),
    @("NORWAY", 20),
    @("HUNGARY", 19),
    @("SOUTH AFRICA", 19),
    @("NIGERIA", 17),
    @("INDONESIA", 16),
    @("LUXEMBOURG", 15),
    @("BANGLADESH", 12),
    @("SLOVENIA", 12),
    @("IRAQ", 11),
    @("CYPRUS", 10),
    @("CROATIA", 9),
    @("KAZAKHSTAN", 9),
    @("MEXICO", 9),
    @("NORTH IRELAND", 8),
    @("UKRAINE", 8),
    @("BELARUS", 7),
    @("SERBIA", 7),
    @("SLOVAKIA", 7),
    @("BULGARIA", 6),
    @("ETHIOPIA", 6),
    @("LITHUANIA", 6),
    @("JORDAN", 5),
    @("LATVIA", 5),
    @("PHILIPPINES", 4),
    @("ALGERIA", 3),
    @("COLOMBIA", 3),
    @("PERU", 3),
    @("SUDAN", 3),
    @("TUNISIA", 3),
    @("ARMENIA", 2),
    @("AZERBAIJAN", 2),
    @("BAHRAIN", 2),
    @("BOSNIA HERCEG", 2),
    @("BRUNEI", 2),
    @("ECUADOR", 2),
    @("GHANA", 2),
    @("MOROCCO", 2),
    @("NORTH MACEDONIA", 2),
    @("SRI LANKA", 2),
    @("BAHAMAS", 1),
    @("BOTSWANA", 1),
    @("CAMBODIA", 1),
    @("ESTONIA", 1),
    @("GEORGIA", 1),
    @("KENYA", 1),
    @("KUWAIT", 1),
    @("LEBANON", 1),
    @("MALTA", 1),
    @("MAURITIUS", 1),
    @("OMAN", 1),
    @("PALESTINE", 1),
    @("PARAGUAY", 1),
    @("VENEZUELA", 1),
    @("YEMEN", 1)
)

$rowCount = $data.Count
$colCount = 2

# Build a 2D array (rows x cols) for bulk assignment
$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $values[$i, 0] = $data[$i][0]
    $values[$i, 1] = $data[$i][1]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1

$range = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 2))
$range.Value = $values
